$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price/Volume columns so numeric-looking strings
# (e.g. "29.355.18", "229.50", "7.410") are preserved exactly as text.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.355.18'
$ws.Range("E2").Value = '  -0.20%  '
$ws.Range("D3").Value = '1.843.23'
$ws.Range("E3").Value = '  -0.45%  '
$ws.Range("D4").Value = '0.9986'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '240.36'
$ws.Range("E5").Value = '  -0.35%  '
$ws.Range("D6").Value = '0.6299'
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("D7").Value = '0.9996'
$ws.Range("E7").Value = '  -0.27%  '
$ws.Range("E8").Value = '  -1.18%  '
$ws.Range("D9").Value = '0.2906'
$ws.Range("E9").Value = '  -0.41%  '
$ws.Range("D10").Value = '24.91'
$ws.Range("E10").Value = '  +1.79%  '
$ws.Range("E11").Value = '  -0.15%  '
$ws.Range("D12").Value = '1.847.21'
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("D13").Value = '4.985'
$ws.Range("E13").Value = '  -0.88%  '
$ws.Range("D14").Value = '0.6781'
$ws.Range("E14").Value = '  -0.56%  '
$ws.Range("D15").Value = '0.00001023'
$ws.Range("E15").Value = '  -1.53%  '
$ws.Range("D16").Value = '82.11'
$ws.Range("E16").Value = '  -1.29%  '
$ws.Range("D17").Value = '6.269'
$ws.Range("E17").Value = '  +2.22%  '
$ws.Range("D18").Value = '29.337.71'
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("D19").Value = '229.50'
$ws.Range("E19").Value = '  +0.04%  '
$ws.Range("D20").Value = '12.32'
$ws.Range("E20").Value = '  -0.36%  '
$ws.Range("D21").Value = '0.9997'
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("D22").Value = '7.410'
$ws.Range("E22").Value = '  -0.77%  '
$ws.Range("E23").Value = '  -0.16%  '
$ws.Range("D24").Value = '158.17'
$ws.Range("E24").Value = '  -0.74%  '
$ws.Range("D25").Value = '8.486'
$ws.Range("E25").Value = '  +0.62%  '
$ws.Range("E26").Value = '  -2.46%  '
$ws.Range("E27").Value = '  -1.25%  '
$ws.Range("D28").Value = '0.06523'
$ws.Range("E28").Value = '  +14.52%  '
$ws.Range("D29").Value = '1.446'
$ws.Range("E29").Value = '  +2.27%  '
$ws.Range("E30").Value = '  +0.67%  '
$ws.Range("D31").Value = '4.071'
$ws.Range("E31").Value = '  -1.67%  '
$ws.Range("D32").Value = '4.055'
$ws.Range("E32").Value = '  +0.18%  '
$ws.Range("E33").Value = '  +1.11%  '
$ws.Range("E34").Value = '  -1.40%  '
$ws.Range("E35").Value = '  -0.17%  '
$ws.Range("E36").Value = '  -0.62%  '
$ws.Range("D37").Value = '0.01855'
$ws.Range("E37").Value = '  +1.36%  '
$ws.Range("D38").Value = '2.816'
$ws.Range("E38").Value = '  -0.94%  '
$ws.Range("D39").Value = '1.247.38'
$ws.Range("E39").Value = '  -0.44%  '
$ws.Range("D40").Value = '6.778'
$ws.Range("E40").Value = '  +4.40%  '
$ws.Range("D41").Value = '0.9313'
$ws.Range("E41").Value = '  +2.53%  '
$ws.Range("D42").Value = '0.9992'
$ws.Range("E42").Value = '  -0.28%  '
$ws.Range("D43").Value = '1.996.89'
$ws.Range("E43").Value = '  -1.02%  '
$ws.Range("D44").Value = '100.79'
$ws.Range("E44").Value = '  -0.72%  '
$ws.Range("D45").Value = '65.55'
$ws.Range("E45").Value = '  -0.74%  '
$ws.Range("E46").Value = '  +2.28%  '
$ws.Range("D47").Value = '7.055'
$ws.Range("E47").Value = '  -0.67%  '
$ws.Range("D48").Value = '1.714'
$ws.Range("E48").Value = '  +2.47%  '
$ws.Range("D49").Value = '9.021'
$ws.Range("E49").Value = '  -0.11%  '
$ws.Range("D50").Value = '0.1146'
$ws.Range("E50").Value = '  -1.44%  '
$ws.Range("D51").Value = '0.3896'
$ws.Range("E51").Value = '  -1.61%  '
